$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data
$ws.Range("D2").Value = "59.652.38"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").Value = "2.664.96"
$ws.Range("E3").Value = "  +2.23%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.89"
$ws.Range("E5").Value = "  -0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.34"
$ws.Range("E6").Value = "  +3.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +0.94%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.79"
$ws.Range("E9").Value = "  +5.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.339"
$ws.Range("E11").Value = "  +1.19%  "

$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").Value = "3.114.92"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").Value = "59.561.57"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.32"
$ws.Range("E15").Value = "  +3.61%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.685.17"
$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.49"
$ws.Range("E18").Value = "  +3.16%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "340.68"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.34"
$ws.Range("E20").Value = "  +2.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.22"

$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.60"
$ws.Range("E23").Value = "  -1.32%  "

$ws.Range("E24").Value = "  +2.47%  "

$ws.Range("E25").Value = "  -0.76%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  +1.28%  "

$ws.Range("E28").Value = "  +1.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("E30").Value = "  -1.44%  "

$ws.Range("E31").Value = "  +0.57%  "

$ws.Range("E32").Value = "  +0.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.78"
$ws.Range("E33").Value = "  +0.91%  "

$ws.Range("E34").Value = "  +0.88%  "

$ws.Range("E35").Value = "  +2.25%  "

$ws.Range("E36").Value = "  +3.19%  "

$ws.Range("E37").Value = "  +0.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.45"
$ws.Range("E38").Value = "  -0.90%  "

$ws.Range("E39").Value = "  +1.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "286.68"
$ws.Range("E40").Value = "  +4.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("E42").Value = "  +1.73%  "

$ws.Range("E43").Value = "  +0.18%  "

$ws.Range("E44").Value = "  +3.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.32"
$ws.Range("E45").Value = "  +4.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0947"
$ws.Range("E46").Value = "  -0.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0228"
$ws.Range("E47").Value = "  +2.14%  "

$ws.Range("D48").Value = "1.967.74"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.41"
$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.16"
$ws.Range("E51").Value = "  -0.38%  "
